$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the description text for the "3. iterace" entry (cell A26)
$ws.Range("A26").Value = "3. iterace - modely balíků, pár dalších scénářů"

# Update hours value for that row (B26): 1 -> 1.5
$ws.Range("B26").Value = 1.5

# Move the selection in the sheet from B27 to A27
$ws.Range("A27").Select()
